# "Generate Report for Handback" — stamp fresh handoff/handback timestamps
# for the 6e192eb7-ea71-4498-b5b3-74c6c8e2fb1e file across the Overview,
# zh-cn and de-de sheets (the 2ade5caa-... row was already handled earlier
# and keeps its existing timestamps).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for row 3 (6e192eb7 file) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Cells.Item(3, 7).Value = "2016-08-22 06:47:18"

# --- zh-cn sheet: row 3 (6e192eb7 file) handoff/handback datetimes ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Cells.Item(3, 8).Value = "2016-08-22 06:47:14"
$zhcn.Cells.Item(3, 11).Value = "2016-08-22 06:47:33"

# --- de-de sheet: row 3 (6e192eb7 file) handoff/handback datetimes ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Cells.Item(3, 8).Value = "2016-08-22 06:47:18"
$dede.Cells.Item(3, 11).Value = "2016-08-22 06:47:40"
